$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 7 description (C7) to append the "save best weights" note
$ws.Range("C7").Value = 'Tested 10, 20, 30 epochs and qualitatively, val_acc does not increase past 10 epochs. Val_acc <= 50%. Considering addding more data aug, manually creating val and train sets, and retraining more layers. Also, considering adding the "save best weights!"'
$ws.Rows.Item(7).RowHeight = 72.5

# New row 8 - save weights functionality added
$ws.Range("A8").Value = "10/12/2019 - 10:05PM"
$ws.Range("B8").Value = "10/12/2019 - 10:42PM"
$ws.Range("C8").Value = 'Added "save weights" upon improvement to functionality. I would add back verbose=1 to see progress bar, or perhaps verbose=2 to see progress bar for at least the epoch. Also, debugged error: all import statements from keras must be specified tensorflow.keras...'
$ws.Range("C8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 72.5

# New row 9 - verbose=2 and save weights confirmed working
$ws.Range("A9").Value = "10/12/2019 - 10:42PM"
$ws.Range("B9").Value = "10/12/2019 - 11:05PM"
$ws.Range("C9").Value = "verbose=2 and save weights implemented correctly!"
$ws.Range("C9").WrapText = $true

$ws.Range("C9").Select()
